# Apply updated nombre_aides (C) and montant_total (D) values
# for the 2020-09-01 data refresh of the Fonds de solidarite volet 1 dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 38788
$ws.Range("D2").Value = 56081959
$ws.Range("C3").Value = 93021
$ws.Range("D3").Value = 136344350
$ws.Range("C4").Value = 31774
$ws.Range("D4").Value = 47052014
$ws.Range("C5").Value = 8920
$ws.Range("D5").Value = 13257397
$ws.Range("C6").Value = 2075
$ws.Range("D6").Value = 3084471
$ws.Range("C7").Value = 171
$ws.Range("D7").Value = 251593
$ws.Range("C12").Value = 42211
$ws.Range("D12").Value = 57227499
$ws.Range("C13").Value = 9899
$ws.Range("D13").Value = 14316446
$ws.Range("C14").Value = 26448
$ws.Range("D14").Value = 38775776
$ws.Range("C15").Value = 8453
$ws.Range("D15").Value = 12544978
$ws.Range("C16").Value = 2206
$ws.Range("D16").Value = 3277539
$ws.Range("C18").Value = 36
$ws.Range("D18").Value = 54000
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10500
$ws.Range("C20").Value = 10411
$ws.Range("D20").Value = 13769778
$ws.Range("C21").Value = 13686
$ws.Range("D21").Value = 19752344
$ws.Range("C22").Value = 32228
$ws.Range("D22").Value = 47285190
$ws.Range("C23").Value = 10405
$ws.Range("D23").Value = 15465048
$ws.Range("C24").Value = 2693
$ws.Range("D24").Value = 4004271
$ws.Range("C25").Value = 537
$ws.Range("D25").Value = 799592
$ws.Range("C27").Value = 11918
$ws.Range("D27").Value = 15908975
$ws.Range("C28").Value = 7863
$ws.Range("D28").Value = 11378638
$ws.Range("C29").Value = 22994
$ws.Range("D29").Value = 33753158
$ws.Range("C30").Value = 7947
$ws.Range("D30").Value = 11819392
$ws.Range("C31").Value = 2008
$ws.Range("D31").Value = 2996251
$ws.Range("C32").Value = 378
$ws.Range("D32").Value = 564415
$ws.Range("C33").Value = 31
$ws.Range("D33").Value = 46393
$ws.Range("C34").Value = 8478
$ws.Range("D34").Value = 11197713
$ws.Range("C35").Value = 3352
$ws.Range("D35").Value = 4841191
$ws.Range("C36").Value = 8021
$ws.Range("D36").Value = 11713926
$ws.Range("C37").Value = 3235
$ws.Range("D37").Value = 4795461
$ws.Range("C38").Value = 838
$ws.Range("D38").Value = 1248223
$ws.Range("C39").Value = 171
$ws.Range("D39").Value = 254186
$ws.Range("C41").Value = 2532
$ws.Range("D41").Value = 3419849
$ws.Range("C42").Value = 17721
$ws.Range("D42").Value = 25623165
$ws.Range("C43").Value = 52198
$ws.Range("D43").Value = 76505887
$ws.Range("C44").Value = 19321
$ws.Range("D44").Value = 28692444
$ws.Range("C45").Value = 5733
$ws.Range("D45").Value = 8534462
$ws.Range("C46").Value = 1252
$ws.Range("D46").Value = 1868545
$ws.Range("C47").Value = 66
$ws.Range("D47").Value = 97068
$ws.Range("C50").Value = 17127
$ws.Range("D50").Value = 22755372
$ws.Range("C51").Value = 2135
$ws.Range("D51").Value = 3098378
$ws.Range("C52").Value = 7204
$ws.Range("D52").Value = 10587251
$ws.Range("C53").Value = 2426
$ws.Range("D53").Value = 3623464
$ws.Range("C54").Value = 770
$ws.Range("D54").Value = 1150305
$ws.Range("C55").Value = 198
$ws.Range("D55").Value = 293226
$ws.Range("C57").Value = 7326
$ws.Range("D57").Value = 10073112
$ws.Range("C58").Value = 1163
$ws.Range("D58").Value = 1961486
$ws.Range("C59").Value = 2850
$ws.Range("D59").Value = 4804854
$ws.Range("C60").Value = 1115
$ws.Range("D60").Value = 1873026
$ws.Range("C61").Value = 387
$ws.Range("D61").Value = 657883
$ws.Range("C62").Value = 126
$ws.Range("D62").Value = 217100
$ws.Range("C63").Value = 26
$ws.Range("D63").Value = 51000
$ws.Range("C64").Value = 1701
$ws.Range("D64").Value = 2660848
$ws.Range("C65").Value = 15787
$ws.Range("D65").Value = 22800689
$ws.Range("C66").Value = 45660
$ws.Range("D66").Value = 66803862
$ws.Range("C67").Value = 15972
$ws.Range("D67").Value = 23730506
$ws.Range("C68").Value = 4647
$ws.Range("D68").Value = 6921309
$ws.Range("C69").Value = 964
$ws.Range("D69").Value = 1434168
$ws.Range("C70").Value = 79
$ws.Range("D70").Value = 115830
$ws.Range("C71").Value = 14
$ws.Range("D71").Value = 19787
$ws.Range("C73").Value = 15400
$ws.Range("D73").Value = 20283351
$ws.Range("C74").Value = 54396
$ws.Range("D74").Value = 79155394
$ws.Range("C75").Value = 152163
$ws.Range("D75").Value = 224157402
$ws.Range("C76").Value = 65658
$ws.Range("D76").Value = 97832049
$ws.Range("C77").Value = 21046
$ws.Range("D77").Value = 31447824
$ws.Range("C78").Value = 5036
$ws.Range("D78").Value = 7522403
$ws.Range("C79").Value = 281
$ws.Range("D79").Value = 416670
$ws.Range("C85").Value = 53418
$ws.Range("D85").Value = 72581016
$ws.Range("C86").Value = 4773
$ws.Range("D86").Value = 6916420
$ws.Range("C87").Value = 11880
$ws.Range("D87").Value = 17448262
$ws.Range("C88").Value = 3966
$ws.Range("D88").Value = 5908958
$ws.Range("C89").Value = 1369
$ws.Range("D89").Value = 2045289
$ws.Range("C90").Value = 295
$ws.Range("D90").Value = 440012
$ws.Range("C92").Value = 6
$ws.Range("D92").Value = 9000
$ws.Range("C93").Value = 5562
$ws.Range("D93").Value = 7475085
$ws.Range("C94").Value = 1653
$ws.Range("D94").Value = 2381699
$ws.Range("C95").Value = 5331
$ws.Range("D95").Value = 7852906
$ws.Range("C96").Value = 1983
$ws.Range("D96").Value = 2952426
$ws.Range("C97").Value = 705
$ws.Range("D97").Value = 1056460
$ws.Range("C98").Value = 197
$ws.Range("D98").Value = 296113
$ws.Range("C101").Value = 3691
$ws.Range("D101").Value = 4888953
$ws.Range("C102").Value = 717
$ws.Range("D102").Value = 1195875
$ws.Range("C103").Value = 440
$ws.Range("D103").Value = 757527
$ws.Range("C104").Value = 160
$ws.Range("D104").Value = 270180
$ws.Range("C105").Value = 52
$ws.Range("D105").Value = 87000
$ws.Range("C106").Value = 27
$ws.Range("D106").Value = 49500
$ws.Range("C107").Value = 11065
$ws.Range("D107").Value = 16053371
$ws.Range("C108").Value = 29743
$ws.Range("D108").Value = 43678970
$ws.Range("C109").Value = 9963
$ws.Range("D109").Value = 14813705
$ws.Range("C110").Value = 2751
$ws.Range("D110").Value = 4101580
$ws.Range("C111").Value = 507
$ws.Range("D111").Value = 755546
$ws.Range("C112").Value = 53
$ws.Range("D112").Value = 79500
$ws.Range("C114").Value = 9982
$ws.Range("D114").Value = 13181024
$ws.Range("C115").Value = 31204
$ws.Range("D115").Value = 44989929
$ws.Range("C116").Value = 67432
$ws.Range("D116").Value = 98665184
$ws.Range("C117").Value = 21735
$ws.Range("D117").Value = 32295788
$ws.Range("C118").Value = 6169
$ws.Range("D118").Value = 9190521
$ws.Range("C119").Value = 1158
$ws.Range("D119").Value = 1730600
$ws.Range("C124").Value = 26326
$ws.Range("D124").Value = 35132698
$ws.Range("C125").Value = 37029
$ws.Range("D125").Value = 53429576
$ws.Range("C126").Value = 78513
$ws.Range("D126").Value = 114790484
$ws.Range("C127").Value = 24279
$ws.Range("D127").Value = 36034041
$ws.Range("C128").Value = 6533
$ws.Range("D128").Value = 9709358
$ws.Range("C129").Value = 1287
$ws.Range("D129").Value = 1913811
$ws.Range("C133").Value = 32467
$ws.Range("D133").Value = 43085105
$ws.Range("C134").Value = 13606
$ws.Range("D134").Value = 19694687
$ws.Range("C135").Value = 32954
$ws.Range("D135").Value = 48393199
$ws.Range("C136").Value = 11672
$ws.Range("D136").Value = 17342087
$ws.Range("C137").Value = 3028
$ws.Range("D137").Value = 4513241
$ws.Range("C139").Value = 36
$ws.Range("D139").Value = 52825
$ws.Range("C141").Value = 11022
$ws.Range("D141").Value = 14689035
$ws.Range("C142").Value = 36128
$ws.Range("D142").Value = 52172139
$ws.Range("C143").Value = 83331
$ws.Range("D143").Value = 122076994
$ws.Range("C144").Value = 24847
$ws.Range("D144").Value = 36911058
$ws.Range("C145").Value = 6532
$ws.Range("D145").Value = 9746496
$ws.Range("C146").Value = 1484
$ws.Range("D146").Value = 2208230
$ws.Range("C147").Value = 84
$ws.Range("D147").Value = 125630
$ws.Range("C149").Value = 29842
$ws.Range("D149").Value = 40229611
